# Rename the embedded logo pictures (Pearson logo x2 in the footers,
# BTec logo in the header) as captured by the commit:
#   - Pearson logo inline pictures: image1.png -> image2.png
#   - BTec logo inline picture:     image2.jpg -> image1.jpg

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Footer 1 (first/"odd" footer) - Pearson Edexcel logo.
$footer1 = $sec.Footers.Item(1)
if ($footer1.Exists -and $footer1.Range.InlineShapes.Count -ge 1) {
    $footer1.Range.InlineShapes.Item(1).Name = "image2.png"
}

# Footer 2 (the other footer) - Pearson Edexcel logo.
$footer2 = $sec.Footers.Item(2)
if ($footer2.Exists -and $footer2.Range.InlineShapes.Count -ge 1) {
    $footer2.Range.InlineShapes.Item(1).Name = "image2.png"
}

# Header 2 - BTec logo.
$header2 = $sec.Headers.Item(2)
if ($header2.Exists -and $header2.Range.InlineShapes.Count -ge 1) {
    $header2.Range.InlineShapes.Item(1).Name = "image1.jpg"
}
